$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36-47 down to 37-48
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with data (mirrors the pattern of the surrounding rows)
$ws.Cells.Item(36, 1).Value = 10
$ws.Cells.Item(36, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(36, 3).Value = "La Araucanía"
$ws.Cells.Item(36, 4).Value = 44609
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 9
$ws.Cells.Item(36, 6).Value = 100114002
$ws.Cells.Item(36, 7).Value = "Camote"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 40
$ws.Cells.Item(36, 11).Value = 18000
$ws.Cells.Item(36, 12).Value = 18000
$ws.Cells.Item(36, 13).Value = 18000
$ws.Cells.Item(36, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(36, 15).Value = "Perú"
$ws.Cells.Item(36, 16).Value = 900
$ws.Cells.Item(36, 17).Value = 20
$ws.Cells.Item(36, 18).Value = "Hortaliza"
